# Edit sheet Card24 by admin
# Remove column O (O1:O12) entirely - the "event" header and its empty
# column cells below it - shrinking the sheet's used range/dimension
# from A1:O12 down to A1:N12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("O1:O12").Clear()
